$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 11112448
$ws.Range("J111").Value = 2415.5
$ws.Range("L111").Value = 7246.5
$ws.Range("N111").Value = -13380.5

$ws.Range("H112").Value = 2371.4546
$ws.Range("J112").Value = 2432
$ws.Range("L112").Value = 7296
$ws.Range("N112").Value = -9512

$ws.Range("H132").Value = 1209.6471
$ws.Range("I132").Value = 1016.42224
$ws.Range("J132").Value = 2658.8333
$ws.Range("K132").Value = 3049.26672
$ws.Range("L132").Value = 7976.499899999999
$ws.Range("M132").Value = -519.2667200000001
$ws.Range("N132").Value = -13036.4999

$ws.Range("H137").Value = 22417.064
$ws.Range("I137").Value = 949.7059
$ws.Range("K137").Value = 2849.1177
$ws.Range("M137").Value = -299.1177000000002

$ws.Range("H138").Value = 1843.7158
$ws.Range("I138").Value = 1557.7258
$ws.Range("J138").Value = 2381.0303
$ws.Range("K138").Value = 4673.1774
$ws.Range("L138").Value = 7143.090899999999
$ws.Range("M138").Value = 466.8226000000004
$ws.Range("N138").Value = -17423.0909

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2645.03
$ws.Range("I32").Value = 2494.6086
$ws.Range("J32").Value = 4374.875
$ws.Range("K32").Value = 2494.6086
$ws.Range("L32").Value = 4374.875
$ws.Range("M32").Value = -2207.6086
$ws.Range("N32").Value = -4948.875

$ws.Range("H61").Value = 4088.7727
$ws.Range("I61").Value = 1310.2667
$ws.Range("K61").Value = 1310.2667
$ws.Range("M61").Value = -1098.2667

$ws.Range("H74").Value = 1462.3043
$ws.Range("J74").Value = 2528.5386
$ws.Range("L74").Value = 2528.5386
$ws.Range("N74").Value = -4276.5386

$ws.Range("H77").Value = 1462.3043
$ws.Range("J77").Value = 2528.5386
$ws.Range("L77").Value = 12642.693
$ws.Range("N77").Value = -21378.693

$ws.Range("H136").Value = 4088.7727
$ws.Range("I136").Value = 1310.2667
$ws.Range("K136").Value = 3930.800099999999
$ws.Range("M136").Value = -1380.800099999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2147.4666
$ws.Range("I20").Value = 1751.5264
$ws.Range("J20").Value = 2831.3635
$ws.Range("K20").Value = 1751.5264
$ws.Range("L20").Value = 2831.3635
$ws.Range("M20").Value = -1504.5264
$ws.Range("N20").Value = -3325.3635

$ws.Range("H99").Value = 1000
$ws.Range("I99").Value = 1000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 498
$ws.Range("N99").ClearContents()

$ws.Range("H134").Value = 4336.3335
$ws.Range("I134").Value = 4384.625
$ws.Range("K134").Value = 13153.875
$ws.Range("M134").Value = -10618.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 10000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 10000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 10000
$ws.Range("N4").Value = -10224
$ws.Range("M4").ClearContents()

$ws.Range("H31").Value = 2281.75
$ws.Range("I31").Value = 1588.8235
$ws.Range("J31").Value = 3352.6365
$ws.Range("K31").Value = 1588.8235
$ws.Range("L31").Value = 3352.6365
$ws.Range("M31").Value = -1293.8235
$ws.Range("N31").Value = -3942.6365

$ws.Range("H34").Value = 2281.75
$ws.Range("I34").Value = 1588.8235
$ws.Range("J34").Value = 3352.6365
$ws.Range("K34").Value = 1588.8235
$ws.Range("L34").Value = 3352.6365
$ws.Range("M34").Value = -1386.8235
$ws.Range("N34").Value = -3756.6365

$ws.Range("H58").Value = 791509.4
$ws.Range("I58").Value = 1115548.2
$ws.Range("J58").Value = 1664.75
$ws.Range("K58").Value = 1115548.2
$ws.Range("L58").Value = 1664.75
$ws.Range("M58").Value = -1115345.2
$ws.Range("N58").Value = -2070.75

$ws.Range("H132").Value = 1676.449
$ws.Range("I132").Value = 1118.8684
$ws.Range("J132").Value = 3602.6365
$ws.Range("K132").Value = 3356.6052
$ws.Range("L132").Value = 10807.9095
$ws.Range("M132").Value = -826.6052
$ws.Range("N132").Value = -15867.9095

$ws.Range("H136").Value = 791509.4
$ws.Range("I136").Value = 1115548.2
$ws.Range("J136").Value = 1664.75
$ws.Range("K136").Value = 3346644.6
$ws.Range("L136").Value = 4994.25
$ws.Range("M136").Value = -3344094.6
$ws.Range("N136").Value = -10094.25

$ws.Range("H141").Value = 80647
$ws.Range("J141").Value = 84970.5
$ws.Range("L141").Value = 84970.5
$ws.Range("N141").Value = -95330.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 91.86667
$ws.Range("I2").Value = 109.90909
$ws.Range("K2").Value = 659.4545400000001
$ws.Range("M2").Value = -546.4545400000001

$ws.Range("H38").Value = 691.8570999999999
$ws.Range("J38").Value = 1334
$ws.Range("L38").Value = 4002
$ws.Range("N38").Value = -4696

$ws.Range("H107").Value = 413.26086
$ws.Range("I107").Value = 312
$ws.Range("J107").Value = 457.5625
$ws.Range("K107").Value = 936
$ws.Range("L107").Value = 1372.6875
$ws.Range("M107").Value = 984
$ws.Range("N107").Value = -5212.6875

$ws.Range("J131").Value = 24433.736
$ws.Range("L131").Value = 73301.208
$ws.Range("N131").Value = -83381.208

$ws.Range("H137").Value = 3661.2856
$ws.Range("J137").Value = 5709.1
$ws.Range("L137").Value = 17127.3
$ws.Range("N137").Value = -27327.3

$ws.Range("H141").Value = 2922.0625
$ws.Range("I141").Value = 2865.6155
$ws.Range("J141").Value = 3166.6667
$ws.Range("K141").Value = 8596.8465
$ws.Range("L141").Value = 9500.000100000001
$ws.Range("M141").Value = -3416.8465
$ws.Range("N141").Value = -19860.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1187.8572
$ws.Range("I97").Value = 1274.2307
$ws.Range("K97").Value = 1274.2307
$ws.Range("M97").Value = -778.2307000000001

$ws.Range("H132").Value = 1070999.5
$ws.Range("I132").Value = 1674373.6
$ws.Range("K132").Value = 5023120.800000001
$ws.Range("M132").Value = -5020590.800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1928.75
$ws.Range("I22").Value = 1750
$ws.Range("K22").Value = 1750
$ws.Range("M22").Value = -1455

$ws.Range("H27").Value = 1928.75
$ws.Range("I27").Value = 1750
$ws.Range("K27").Value = 1750
$ws.Range("M27").Value = -1643

$ws.Range("H97").Value = 14999.5
$ws.Range("J97").Value = 14999.5
$ws.Range("L97").Value = 14999.5
$ws.Range("N97").Value = -16981.5

$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws.Range("H136").Value = 1682.9615
$ws.Range("I136").Value = 1192.8029
$ws.Range("J136").Value = 6654.5713
$ws.Range("K136").Value = 3578.4087
$ws.Range("L136").Value = 19963.7139
$ws.Range("M136").Value = -1028.4087
$ws.Range("N136").Value = -25063.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8050.4287
$ws.Range("I62").Value = 8110.6
$ws.Range("K62").Value = 8110.6
$ws.Range("M62").Value = -7486.6

$ws.Range("H65").Value = 8050.4287
$ws.Range("I65").Value = 8110.6
$ws.Range("K65").Value = 40553
$ws.Range("M65").Value = -37433

$ws.Range("H122").Value = 69507.22
$ws.Range("I122").Value = 99191.31
$ws.Range("J122").Value = 1657.8572
$ws.Range("K122").Value = 297573.93
$ws.Range("L122").Value = 4973.571599999999
$ws.Range("M122").Value = -295123.93
$ws.Range("N122").Value = -9873.571599999999

$ws.Range("H126").Value = 6275.423
$ws.Range("I126").Value = 6761.1577
$ws.Range("K126").Value = 20283.4731
$ws.Range("M126").Value = -17813.4731

$ws.Range("H132").Value = 1285.4
$ws.Range("I132").Value = 879.3617
$ws.Range("K132").Value = 2638.0851
$ws.Range("M132").Value = -108.0851000000002

$ws.Range("H136").Value = 25256464
$ws.Range("I136").Value = 37041480
$ws.Range("K136").Value = 111124440
$ws.Range("M136").Value = -111121890
